{"js": "// Replace the three-digit-by-one-digit division expressions throughout\n// the document body with the new set of expressions, matching the\n// unified diff (each cell's \"NNN\u00f7N=\" text is replaced by a new one).\nconst replacements = [\n  [\"606\u00f79=\", \"308\u00f75=\"],\n  [\"657\u00f78=\", \"820\u00f76=\"],\n  [\"735\u00f73=\", \"731\u00f74=\"],\n  [\"639\u00f77=\", \"536\u00f74=\"],\n  [\"532\u00f79=\", \"734\u00f76=\"],\n  [\"201\u00f75=\", \"794\u00f79=\"],\n  [\"916\u00f72=\", \"550\u00f78=\"],\n  [\"322\u00f78=\", \"382\u00f73=\"],\n  [\"602\u00f78=\", \"463\u00f74=\"],\n  [\"209\u00f75=\", \"370\u00f79=\"],\n  [\"649\u00f75=\", \"711\u00f75=\"],\n  [\"503\u00f75=\", \"295\u00f77=\"],\n  [\"671\u00f79=\", \"378\u00f74=\"],\n  [\"554\u00f72=\", \"613\u00f78=\"],\n  [\"140\u00f79=\", \"126\u00f76=\"],\n  [\"371\u00f72=\", \"758\u00f77=\"],\n  [\"577\u00f74=\", \"764\u00f72=\"],\n  [\"748\u00f76=\", \"118\u00f74=\"],\n  [\"699\u00f76=\", \"645\u00f76=\"],\n  [\"686\u00f78=\", \"789\u00f76=\"],\n  [\"829\u00f77=\", \"465\u00f78=\"],\n  [\"922\u00f79=\", \"636\u00f77=\"],\n  [\"474\u00f76=\", \"491\u00f77=\"],\n  [\"594\u00f76=\", \"826\u00f73=\"],\n  [\"266\u00f75=\", \"568\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit division expressions throughout\n# the document body with the new set of expressions, matching the\n# unified diff (each cell's \"NNN\u00f7N=\" text is replaced by a new one).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"606\u00f79=\", \"308\u00f75=\"),\n    @(\"657\u00f78=\", \"820\u00f76=\"),\n    @(\"735\u00f73=\", \"731\u00f74=\"),\n    @(\"639\u00f77=\", \"536\u00f74=\"),\n    @(\"532\u00f79=\", \"734\u00f76=\"),\n    @(\"201\u00f75=\", \"794\u00f79=\"),\n    @(\"916\u00f72=\", \"550\u00f78=\"),\n    @(\"322\u00f78=\", \"382\u00f73=\"),\n    @(\"602\u00f78=\", \"463\u00f74=\"),\n    @(\"209\u00f75=\", \"370\u00f79=\"),\n    @(\"649\u00f75=\", \"711\u00f75=\"),\n    @(\"503\u00f75=\", \"295\u00f77=\"),\n    @(\"671\u00f79=\", \"378\u00f74=\"),\n    @(\"554\u00f72=\", \"613\u00f78=\"),\n    @(\"140\u00f79=\", \"126\u00f76=\"),\n    @(\"371\u00f72=\", \"758\u00f77=\"),\n    @(\"577\u00f74=\", \"764\u00f72=\"),\n    @(\"748\u00f76=\", \"118\u00f74=\"),\n    @(\"699\u00f76=\", \"645\u00f76=\"),\n    @(\"686\u00f78=\", \"789\u00f76=\"),\n    @(\"829\u00f77=\", \"465\u00f78=\"),\n    @(\"922\u00f79=\", \"636\u00f77=\"),\n    @(\"474\u00f76=\", \"491\u00f77=\"),\n    @(\"594\u00f76=\", \"826\u00f73=\"),\n    @(\"266\u00f75=\", \"568\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
